$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update status text "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Update timestamps
$wsOverview.Range("G2").Value = "2016-08-14 01:18:50"
$wsDeDe.Range("H2").Value = "2016-08-14 01:18:50"
$wsZhCn.Range("H2").Value = "2016-08-14 01:18:42"

# Widen columns to fit new, longer text ("Ready for handoff").
# Target stored width is ~17.216; the host's ColumnWidth setter only
# persists widths on a whole-pixel (1/6 character unit) grid, so 16.25
# is the input that lands on the closest achievable grid point (17.1667).
$wsOverview.Range("E:E").ColumnWidth = 16.25
$wsOverview.Range("F:F").ColumnWidth = 16.25
$wsZhCn.Range("C:C").ColumnWidth = 16.25
$wsDeDe.Range("C:C").ColumnWidth = 16.25
